# Bugfixes in VPC plots and response on last review (#20)
#
# - "outputxls" (C14) no longer hardcodes "Workflow.xlsx" -> clear it
# - "calculatePKParameterFh" (C16) now points to a custom function handle
# - "TasksimulatePopulation" (C26) switched from true (1) to false (0)
# - Selection / scroll position on the Workflow sheet moved up (A10 / C18)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workflow")

# Clear the old "Workflow.xlsx" value for outputxls
$ws.Range("C14").Value = ""

# Set the new function handle for calculatePKParameterFh
$ws.Range("C16").Value = "myCalculatePKParameterForApplicationProtocol"

# TasksimulatePopulation flips from 1 (true) to 0 (false)
$ws.Range("C26").Value = 0

# Update the active view / selection for the Workflow sheet
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C18").Select()
